$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two disciplinary records that were merged/cleared out entirely
# (the rows for "Nguyễn Văn A" and "Nguyễn Thị E"). Deleting row 4 first keeps
# row 3's index valid for the second delete.
$ws.Rows.Item(4).EntireRow.Delete()
$ws.Rows.Item(3).EntireRow.Delete()

# Row 2 ("Trần Thị Long Lanh Sương Sớm Mai"): keep only the "Đi trễ" violation
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = "Đi trễ 12/11/2023 (admin, 20)"

# Row 3 ("Nguyễn Lê Phi Long"): drop the trailing semicolon
$ws.Cells.Item(3, 5).Value = "Đi trễ 12/11/2023 (admin, 21)"

# Row 4 ("Phan Hoài Linh"): keep only the three 12/11 violations, one per line
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = "Ăn quà sai quy định 12/11/2023 (admin, 17)`nĂn quà sai quy định 12/11/2023 (admin, 18)`nĂn quà sai quy định 12/11/2023 (admin, 19)"

# Apply the header's cell style (wrap text) to all the data rows
$ws.Range("A1:E1").Copy()
$ws.Range("A2:E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Narrow column E now that the text is much shorter
$ws.Columns.Item(5).ColumnWidth = 50.559

# Match the new selection/active cell recorded in the saved file
$ws.Range("E4").Select()
